$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 177. This shifts the existing rows 177-225
# down to 178-226 (carrying all their data/formatting with them), and
# leaves a blank row 177 for the new record.
$ws.Rows.Item(177).Insert()

# Populate the new row 177 with the new data record. The columns that
# stay identical to the (now shifted) neighboring record are filled in
# explicitly as well, matching the target state.
$ws.Range("A177").Value = 9
$ws.Range("B177").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C177").Value = 'Metropolitana'
$ws.Range("D177").Value = 44588
$ws.Range("E177").Value = 13
$ws.Range("F177").Value = 100112043
$ws.Range("G177").Value = 'Pepino ensalada'
$ws.Range("H177").Value = 'Sin especificar'
$ws.Range("I177").Value = 'Primera'
$ws.Range("J177").Value = 106
$ws.Range("K177").Value = 11000
$ws.Range("L177").Value = 12000
$ws.Range("M177").Value = 11500
$ws.Range("N177").Value = '$/caja 60 unidades'
$ws.Range("O177").Value = 'Región de Arica y Parinacota'
$ws.Range("P177").Value = 192
$ws.Range("Q177").Value = 60
$ws.Range("R177").Value = 'Hortaliza'
